$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-LatticeCell($table, $row, $col, $title, $line2, $d1, $d2) {
    $cell = $table.Cell($row, $col)
    $xmlSnippet = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>' + $title + '</w:t><w:br/><w:t xml:space="preserve">' + $line2 + '</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>' + $d1 + '</w:t><w:br/><w:t>' + $d2 + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $cell.Range.InsertXML($xmlSnippet)
}

Set-LatticeCell $t 1 1 "23 x 89" "  8    9" "2|    |" "3|    |"
Set-LatticeCell $t 1 2 "21 x 50" "  5    0" "2|    |" "1|    |"
Set-LatticeCell $t 1 3 "16 x 33" "  3    3" "1|    |" "6|    |"
Set-LatticeCell $t 2 1 "38 x 21" "  2    1" "3|    |" "8|    |"
Set-LatticeCell $t 2 2 "22 x 58" "  5    8" "2|    |" "2|    |"
Set-LatticeCell $t 2 3 "50 x 72" "  7    2" "5|    |" "0|    |"
Set-LatticeCell $t 3 1 "28 x 74" "  7    4" "2|    |" "8|    |"
Set-LatticeCell $t 3 2 "46 x 64" "  6    4" "4|    |" "6|    |"
Set-LatticeCell $t 3 3 "67 x 33" "  3    3" "6|    |" "7|    |"
Set-LatticeCell $t 4 1 "50 x 64" "  6    4" "5|    |" "0|    |"
Set-LatticeCell $t 4 2 "42 x 59" "  5    9" "4|    |" "2|    |"
Set-LatticeCell $t 4 3 "32 x 39" "  3    9" "3|    |" "2|    |"
Set-LatticeCell $t 5 1 "77 x 55" "  5    5" "7|    |" "7|    |"
Set-LatticeCell $t 5 2 "40 x 83" "  8    3" "4|    |" "0|    |"
Set-LatticeCell $t 5 3 "56 x 13" "  1    3" "5|    |" "6|    |"
